$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: remove the old 3rd data row entirely ---
$ws1.Rows.Item(3).Delete()

# --- Sheet1: update header row (A1:F1); G1 is set later (it is the last
# new shared string the original author typed) ---
$ws1.Range("A1").Value = "ID sản phẩm"
$ws1.Range("B1").Value = "Số lô"
$ws1.Range("C1").Value = "Số lượng"
$ws1.Range("D1").Value = "Giá gốc/sản phẩm"
$ws1.Range("E1").Value = "Ngày nhập"
$ws1.Range("F1").Value = "ID nhà cung cấp"

# --- Sheet1: update row 2 data values ---
$ws1.Range("A2").Value = 144
$ws1.Range("C2").Value = 1440
$ws1.Range("D2").Value = 2000
$ws1.Range("F2").Value = 2

# --- Sheet1: new formula cells (date / date-time) ---
$ws1.Range("E2").Formula = "=TODAY()"
$ws1.Range("E2").NumberFormat = "mm-dd-yy"
$ws1.Range("G2").Formula = "=NOW()"
$ws1.Range("G2").NumberFormat = "m/d/yy h:mm"

# --- Sheet1: column widths ---
$ws1.Columns.Item(1).ColumnWidth = 13.42578125
$ws1.Columns.Item(2).ColumnWidth = 15.140625
$ws1.Columns.Item(3).ColumnWidth = 12.28515625
$ws1.Columns.Item(4).ColumnWidth = 18.42578125
$ws1.Columns.Item(5).ColumnWidth = 13.42578125
$ws1.Columns.Item(6).ColumnWidth = 16.85546875
$ws1.Columns.Item(7).ColumnWidth = 16.140625
$ws1.Columns.Item(8).ColumnWidth = 13.7109375

# --- Sheet1: page orientation ---
$ws1.PageSetup.Orientation = 1

# --- Add Sheet2 after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Sheet2: header row ---
$ws2.Range("A1").Value = "productID"
$ws2.Range("B1").Value = "productName"
$ws2.Range("C1").Value = "productPrice"
$ws2.Range("D1").Value = "brand"
$ws2.Range("E1").Value = "color"
$ws2.Range("F1").Value = "material"
$ws2.Range("G1").Value = "weight"
$ws2.Range("H1").Value = "dimensions"

# --- Sheet2: row 2 data (order matters for shared-string layout) ---
$ws2.Range("A2").Value = 144
$ws2.Range("C2").Value = 8000
$ws2.Range("D2").Value = "CapyTech"
$ws2.Range("E2").Value = "Cam"
$ws2.Range("F2").Value = "Kim loại"
$ws2.Range("G2").Value = 0.9
$ws2.Range("H2").Value = "A4"
$ws2.Range("B2").Value = "Giấy nhớ 5x"

# --- Sheet1: final header cell (typed last by the original author) ---
$ws1.Range("G1").Value = "Ngày giờ tạo"

# --- Sheet2: column widths ---
$ws2.Columns.Item(1).ColumnWidth = 15.85546875
$ws2.Columns.Item(2).ColumnWidth = 13.7109375
$ws2.Columns.Item(3).ColumnWidth = 14.5703125
$ws2.Columns.Item(4).ColumnWidth = 11.7109375

# --- Selections ---
$ws1.Range("F5").Select()
$ws2.Range("F7").Select()
$ws1.Activate()
$ws1.Range("F5").Select()

# --- Workbook window position ---
$excel.ActiveWindow.Left = 4455
$excel.ActiveWindow.Top = 1305
